$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the hyperlink that currently lives on B2 (it's about to become a
# plain price cell instead of the Ebay URL cell).
# ---------------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Row 1 header relabelling: B1 becomes "Price", F1 becomes "Ebay URL".
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Price"
$ws.Range("F1").Value = "Ebay URL"

# ---------------------------------------------------------------------------
# Row 2: B2 becomes the plain price text (style back to Normal so the old
# hyperlink formatting doesn't linger), F2 gets the Ebay URL that used to
# live in B2, formatted + linked as a hyperlink.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "US `$348.00"
$ws.Range("B2").Style = "Normal"

$ws.Range("F2").Value = "https://www.ebay.com/itm/VisionTek-RX-480-8GB-GDDR5-Overclocked-Edition-Rear-Blower-4M-3x-DP-HDMI/114683471280?hash=item1ab3aae5b0:g:aqoAAOSwK1tgKFVK"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.ebay.com/itm/VisionTek-RX-480-8GB-GDDR5-Overclocked-Edition-Rear-Blower-4M-3x-DP-HDMI/114683471280?hash=item1ab3aae5b0:g:aqoAAOSwK1tgKFVK")
$ws.Range("F2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Row 3 (new): second watched item's Ebay URL.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "https://www.ebay.com/itm/nvidia-rtx-3060-ti-founders-edition/203288698285?hash=item2f54f35dad:g:gLkAAOSwraBgNu-p"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.ebay.com/itm/nvidia-rtx-3060-ti-founders-edition/203288698285?hash=item2f54f35dad:g:gLkAAOSwraBgNu-p")
$ws.Range("F3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Row 4 (new): third watched item's Ebay URL, with an explicit display text
# (same as the address) to mirror the source diff's "display" attribute.
# ---------------------------------------------------------------------------
$longUrl = "https://www.ebay.com/itm/NVIDIA-GeForce-RTX-2080-Ti-Cyberpunk-2077-Edition-Front-Back-Plate-Only-USED/174648293981?_trkparms=aid%3D1110006%26algo%3DHOMESPLICE.SIM%26ao%3D1%26asc%3D230925%26meid%3D43882e1369e9447d8cea1e7aebe460c6%26pid%3D101195%26rk%3D2%26rkt%3D12%26mehot%3Dpf%26sd%3D203288698285%26itm%3D174648293981%26pmt%3D1%26noa%3D0%26pg%3D2047675%26algv%3DSimplAMLv9PairwiseUnbiasedWeb%26brand%3DNVIDIA&_trksid=p2047675.c101195.m1851"
$ws.Range("F4").Value = $longUrl
$ws.Hyperlinks.Add($ws.Range("F4"), $longUrl, [System.Type]::Missing, $longUrl)
$ws.Range("F4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column layout: drop the custom widths on B & C (back to sheet default),
# give F the width that used to belong to the URL column.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 8
$ws.Columns("C").ColumnWidth = 8
$ws.Columns("F").ColumnWidth = 35.140625

# ---------------------------------------------------------------------------
# Selection / view bits to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("C11").Select()

Write-Output "done"
